$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the statement labels used as "Parent Public ID" values
# (update every cell that shares the old label so all references move together)
$ws.Range("A2").Value = "statement-01"
$ws.Range("A3").Value = "statement-01"
$ws.Range("A4").Value = "statement-01"
$ws.Range("A5").Value = "statement-02"
$ws.Range("A6").Value = "statement-02"

# Clear the auto-generated "Public ID" values for the measure rows,
# keeping their existing formatting/style intact.
$ws.Range("B2:B6").ClearContents()
